# Split the single merged field-open run "{m" + ":userdoc 'zone1'}" in the
# second paragraph into four separate runs: "{", "m", ":userdoc 'zone1'", "}"
# (mirrors TokenIteratorFieldRewriterSplit's tokenisation of the m:userdoc
# field marker).

$d = $word.ActiveDocument

# Locate the paragraph that currently reads "{m:userdoc 'zone1'}".
# (Paragraph.Range.Text includes the trailing paragraph mark, so trim it
# before comparing.)
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "{m:userdoc 'zone1'}") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find paragraph containing {m:userdoc 'zone1'}"
}

$r = $target.Range
$start = $r.Start
$end = $r.End

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:r><w:t>{</w:t></w:r>' +
       '<w:r><w:t>m</w:t></w:r>' +
       '<w:r><w:t>:userdoc ''zone1''</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
       '</w:p>'

$fullRange = $d.Range($start, $end)
$fullRange.InsertXML($xml)
